# Standardize FK column names (and a couple of related labels) with underscores.
# Only the Field_Name column (D) is updated to the canonical underscore format,
# except for row 45 where both Field_Name (D) and Field_Label (E) are corrected
# from the stale "Master[EBITA]" value to "EBITA".
#
# NOTE on ordering: new shared-string entries are appended to the workbook's
# shared-string table in the order the cells are written, so the write order
# below (Standard_costing / Account_Name first, then the FK_* underscore
# names) intentionally matches the order the strings appear appended at the
# end of xl/sharedStrings.xml in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D32").Value = "Standard_costing"
$ws.Range("D41").Value = "Account_Name"

$ws.Range("D9").Value  = "FK_location_field_model"
$ws.Range("D10").Value = "FK_location_field_model_id"
$ws.Range("D11").Value = "FK_location_record_Id"
$ws.Range("D12").Value = "Qdrant_ID_for_FK"

$ws.Range("D45").Value = "EBITA"
$ws.Range("E45").Value = "EBITA"

# Match the author's final cell selection recorded in the workbook view.
[void]$ws.Range("E19").Select()
